# Add/update metadata report for Akurana
# The new row (row 4) duplicates row 3's data (Dec 2024 / Akurana entry with
# no data availability), so copy row 3 verbatim into row 4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$sourceRow = $ws.Range("A3:AO3")
$targetRow = $ws.Range("A4:AO4")

$targetRow.Value2 = $sourceRow.Value2
